$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 6 new Hungarian participant rows after existing row 30 ---
$hungarianRows = @(
    @{ Pid = 'P3SEIZnX6k3'; L1 = 'Hungarian'; Data = '[{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"dør","word2":"dør"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"ful","word2":"ful"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"sinde","word2":"sende"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ful","word2":"fyl"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"fuld","word2":"fugl"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"før","word2":"fær"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"kat","word2":"kæt"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"kat","word2":"kat"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"rød","word2":"ryd"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ven","word2":"pen"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"ben","word2":"ben"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"gul","word2":"guld"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"rør","word2":"rær"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"grå","word2":"grå"},{"participantResponse":"different","correctResponse":"same","isCorrect":false,"word1":"syr","word2":"syr"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rød","word2":"lød"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ben","word2":"bøn"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"lys","word2":"los"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"sø","word2":"su"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"lys","word2":"lys"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"løg","word2":"ly"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"fisk","word2":"fæsk"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"syn","word2":"søn"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"lidt","word2":"let"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"rød","word2":"rød"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"fuld","word2":"fuld"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"set","word2":"sæt"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"mølle","word2":"mulle"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"røre","word2":"røre"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"skæl","word2":"skæl"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"nø","word2":"nø"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"møl","word2":"møl"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"pil","word2":"bil"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"syd","word2":"sød"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"væg","word2":"vægt"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"tat","word2":"tæt"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"mæt","word2":"mæt"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"nø","word2":"nu"}]' },
    @{ Pid = 'rXaJ2sWSnuk'; L1 = 'Hungarian'; Data = '[{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"nø","word2":"nu"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"fuld","word2":"fugl"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"syn","word2":"søn"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"syd","word2":"sød"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"tat","word2":"tæt"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"skæl","word2":"skæl"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"set","word2":"sæt"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"nø","word2":"nø"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"møl","word2":"møl"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"fuld","word2":"fuld"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"mølle","word2":"mulle"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"lidt","word2":"let"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"fisk","word2":"fæsk"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"sø","word2":"su"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"gul","word2":"guld"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"grå","word2":"grå"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ben","word2":"bøn"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"røre","word2":"røre"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rød","word2":"ryd"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"kat","word2":"kæt"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"lys","word2":"los"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"før","word2":"fær"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"rød","word2":"rød"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"dør","word2":"dør"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ven","word2":"pen"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"kat","word2":"kat"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ful","word2":"fyl"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"sinde","word2":"sende"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rør","word2":"rær"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rød","word2":"lød"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"lys","word2":"lys"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"ful","word2":"ful"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"løg","word2":"ly"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"syr","word2":"syr"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"væg","word2":"vægt"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"mæt","word2":"mæt"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"ben","word2":"ben"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"pil","word2":"bil"}]' },
    @{ Pid = 'GSiXF9sB7Lu'; L1 = 'Hungarian'; Data = '[{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"ben","word2":"ben"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"tat","word2":"tæt"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ben","word2":"bøn"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"før","word2":"fær"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"fuld","word2":"fugl"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ven","word2":"pen"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"gul","word2":"guld"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"kat","word2":"kæt"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"fuld","word2":"fuld"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"skæl","word2":"skæl"},{"participantResponse":"different","correctResponse":"same","isCorrect":false,"word1":"ful","word2":"ful"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"pil","word2":"bil"},{"participantResponse":"different","correctResponse":"same","isCorrect":false,"word1":"nø","word2":"nø"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"nø","word2":"nu"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"set","word2":"sæt"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"dør","word2":"dør"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"røre","word2":"røre"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"lidt","word2":"let"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"mæt","word2":"mæt"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"løg","word2":"ly"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"lys","word2":"los"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"lys","word2":"lys"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"grå","word2":"grå"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rød","word2":"ryd"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"mølle","word2":"mulle"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rør","word2":"rær"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"rød","word2":"rød"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rød","word2":"lød"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ful","word2":"fyl"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"syd","word2":"sød"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"sinde","word2":"sende"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"syn","word2":"søn"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"kat","word2":"kat"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"fisk","word2":"fæsk"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"sø","word2":"su"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"væg","word2":"vægt"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"møl","word2":"møl"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"syr","word2":"syr"}]' },
    @{ Pid = 'X3YZI52B3Qe'; L1 = 'Hungarian'; Data = '[{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"lys","word2":"los"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"mølle","word2":"mulle"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"fuld","word2":"fugl"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"syd","word2":"sød"},{"participantResponse":"different","correctResponse":"same","isCorrect":false,"word1":"lys","word2":"lys"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"syn","word2":"søn"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"før","word2":"fær"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"nø","word2":"nu"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"kat","word2":"kat"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"rød","word2":"rød"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"ben","word2":"ben"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"fisk","word2":"fæsk"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"fuld","word2":"fuld"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"væg","word2":"vægt"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"gul","word2":"guld"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"kat","word2":"kæt"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rør","word2":"rær"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"løg","word2":"ly"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rød","word2":"lød"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"sinde","word2":"sende"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"sø","word2":"su"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"set","word2":"sæt"},{"participantResponse":"different","correctResponse":"same","isCorrect":false,"word1":"nø","word2":"nø"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"mæt","word2":"mæt"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"pil","word2":"bil"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"tat","word2":"tæt"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"lidt","word2":"let"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"møl","word2":"møl"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ful","word2":"fyl"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"skæl","word2":"skæl"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ven","word2":"pen"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ben","word2":"bøn"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"ful","word2":"ful"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"syr","word2":"syr"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"grå","word2":"grå"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"røre","word2":"røre"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"rød","word2":"ryd"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"dør","word2":"dør"}]' },
    @{ Pid = 'ndIPuTc6WNs'; L1 = 'Hungarian'; Data = '[{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"lidt","word2":"let"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"fuld","word2":"fugl"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ven","word2":"pen"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"røre","word2":"røre"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"rød","word2":"rød"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"dør","word2":"dør"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"grå","word2":"grå"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ben","word2":"bøn"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"tat","word2":"tæt"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"lys","word2":"los"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"sø","word2":"su"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"syd","word2":"sød"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"syr","word2":"syr"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"skæl","word2":"skæl"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"før","word2":"fær"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"nø","word2":"nø"},{"participantResponse":"different","correctResponse":"same","isCorrect":false,"word1":"mæt","word2":"mæt"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"rød","word2":"ryd"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"gul","word2":"guld"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rød","word2":"lød"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"mølle","word2":"mulle"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"set","word2":"sæt"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"møl","word2":"møl"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"ful","word2":"ful"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"fisk","word2":"fæsk"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"ben","word2":"ben"},{"participantResponse":"different","correctResponse":"same","isCorrect":false,"word1":"lys","word2":"lys"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"væg","word2":"vægt"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"kat","word2":"kæt"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"fuld","word2":"fuld"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"syn","word2":"søn"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ful","word2":"fyl"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rør","word2":"rær"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"sinde","word2":"sende"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"løg","word2":"ly"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"kat","word2":"kat"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"pil","word2":"bil"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"nø","word2":"nu"}]' },
    @{ Pid = 'UpHbfcJxkeJ'; L1 = 'Hungarian'; Data = '[{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"mæt","word2":"mæt"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rød","word2":"ryd"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"pil","word2":"bil"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"mølle","word2":"mulle"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"dør","word2":"dør"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ful","word2":"fyl"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"fuld","word2":"fuld"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"sinde","word2":"sende"},{"participantResponse":"different","correctResponse":"same","isCorrect":false,"word1":"grå","word2":"grå"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"væg","word2":"vægt"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"set","word2":"sæt"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"løg","word2":"ly"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"fisk","word2":"fæsk"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"før","word2":"fær"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"sø","word2":"su"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"kat","word2":"kat"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"nø","word2":"nu"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"lidt","word2":"let"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"tat","word2":"tæt"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"lys","word2":"lys"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"kat","word2":"kæt"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"fuld","word2":"fugl"},{"participantResponse":"different","correctResponse":"same","isCorrect":false,"word1":"ben","word2":"ben"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ben","word2":"bøn"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ven","word2":"pen"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"nø","word2":"nø"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"syn","word2":"søn"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"gul","word2":"guld"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"lys","word2":"los"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"rør","word2":"rær"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"møl","word2":"møl"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"skæl","word2":"skæl"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"syd","word2":"sød"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"ful","word2":"ful"},{"participantResponse":"different","correctResponse":"same","isCorrect":false,"word1":"rød","word2":"rød"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"røre","word2":"røre"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rød","word2":"lød"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"syr","word2":"syr"}]' }
)

$insertAt = 31
foreach ($row in $hungarianRows) {
    $ws.Rows.Item($insertAt).Insert()
    $ws.Range("A" + $insertAt).Value = $row.Pid
    $ws.Range("B" + $insertAt).Value = $row.L1
    $ws.Range("C" + $insertAt).Value = $row.Data
    $insertAt = $insertAt + 1
}

# --- Append 6 new Spanish participant rows at the end of the sheet ---
$spanishRows = @(
    @{ Pid = 'XcVhPNabrOL'; L1 = 'Spanish'; Data = '[{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"fuld","word2":"fugl"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rød","word2":"lød"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ven","word2":"pen"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"sinde","word2":"sende"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"møl","word2":"møl"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"væg","word2":"vægt"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"sø","word2":"su"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"set","word2":"sæt"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"syd","word2":"sød"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ben","word2":"bøn"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"kat","word2":"kat"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"fuld","word2":"fuld"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"grå","word2":"grå"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"dør","word2":"dør"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"lys","word2":"lys"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"løg","word2":"ly"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"lidt","word2":"let"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"tat","word2":"tæt"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"mæt","word2":"mæt"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ful","word2":"fyl"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rød","word2":"ryd"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"rød","word2":"rød"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"gul","word2":"guld"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"før","word2":"fær"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"mølle","word2":"mulle"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"nø","word2":"nø"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"syn","word2":"søn"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"ben","word2":"ben"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"fisk","word2":"fæsk"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"kat","word2":"kæt"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"skæl","word2":"skæl"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"pil","word2":"bil"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"lys","word2":"los"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"syr","word2":"syr"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"ful","word2":"ful"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"røre","word2":"røre"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"nø","word2":"nu"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rør","word2":"rær"}]' },
    @{ Pid = 'XWLqjVHufIJ'; L1 = 'Spanish'; Data = '[{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ful","word2":"fyl"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"før","word2":"fær"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rør","word2":"rær"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"fuld","word2":"fuld"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"fuld","word2":"fugl"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"møl","word2":"møl"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"mæt","word2":"mæt"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"sinde","word2":"sende"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"lidt","word2":"let"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"pil","word2":"bil"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"gul","word2":"guld"},{"participantResponse":"different","correctResponse":"same","isCorrect":false,"word1":"nø","word2":"nø"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ven","word2":"pen"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"kat","word2":"kæt"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"ben","word2":"ben"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"røre","word2":"røre"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"lys","word2":"los"},{"participantResponse":"different","correctResponse":"same","isCorrect":false,"word1":"rød","word2":"rød"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"fisk","word2":"fæsk"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rød","word2":"ryd"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"set","word2":"sæt"},{"participantResponse":"different","correctResponse":"same","isCorrect":false,"word1":"ful","word2":"ful"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ben","word2":"bøn"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"syr","word2":"syr"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rød","word2":"lød"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"syd","word2":"sød"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"sø","word2":"su"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"dør","word2":"dør"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"løg","word2":"ly"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"væg","word2":"vægt"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"nø","word2":"nu"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"mølle","word2":"mulle"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"grå","word2":"grå"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"syn","word2":"søn"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"skæl","word2":"skæl"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"kat","word2":"kat"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"tat","word2":"tæt"},{"participantResponse":"different","correctResponse":"same","isCorrect":false,"word1":"lys","word2":"lys"}]' },
    @{ Pid = '2SjTHoDevt6'; L1 = 'Spanish'; Data = '[{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"lys","word2":"los"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"skæl","word2":"skæl"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ful","word2":"fyl"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"gul","word2":"guld"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"fuld","word2":"fuld"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"set","word2":"sæt"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"før","word2":"fær"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"møl","word2":"møl"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"grå","word2":"grå"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rør","word2":"rær"},{"participantResponse":"different","correctResponse":"same","isCorrect":false,"word1":"ben","word2":"ben"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"røre","word2":"røre"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"syr","word2":"syr"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"nø","word2":"nu"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"syd","word2":"sød"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"lys","word2":"lys"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rød","word2":"lød"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"fuld","word2":"fugl"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ben","word2":"bøn"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"lidt","word2":"let"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"fisk","word2":"fæsk"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"løg","word2":"ly"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"dør","word2":"dør"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"mæt","word2":"mæt"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"mølle","word2":"mulle"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ven","word2":"pen"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"nø","word2":"nø"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"ful","word2":"ful"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"kat","word2":"kæt"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"rød","word2":"rød"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"sinde","word2":"sende"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"kat","word2":"kat"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rød","word2":"ryd"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"syn","word2":"søn"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"pil","word2":"bil"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"tat","word2":"tæt"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"væg","word2":"vægt"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"sø","word2":"su"}]' },
    @{ Pid = 'zchpgnvwRAe'; L1 = 'Spanish'; Data = '[{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"fisk","word2":"fæsk"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"syd","word2":"sød"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"lidt","word2":"let"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"mølle","word2":"mulle"},{"participantResponse":"different","correctResponse":"same","isCorrect":false,"word1":"mæt","word2":"mæt"},{"participantResponse":"different","correctResponse":"same","isCorrect":false,"word1":"ben","word2":"ben"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"nø","word2":"nu"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"grå","word2":"grå"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"før","word2":"fær"},{"participantResponse":"different","correctResponse":"same","isCorrect":false,"word1":"røre","word2":"røre"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"tat","word2":"tæt"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"dør","word2":"dør"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"kat","word2":"kæt"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"nø","word2":"nø"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"rød","word2":"rød"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rød","word2":"ryd"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"løg","word2":"ly"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"lys","word2":"lys"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"fuld","word2":"fugl"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"lys","word2":"los"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ben","word2":"bøn"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"sinde","word2":"sende"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rød","word2":"lød"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ven","word2":"pen"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"set","word2":"sæt"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"møl","word2":"møl"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"kat","word2":"kat"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"sø","word2":"su"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"skæl","word2":"skæl"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"pil","word2":"bil"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"ful","word2":"ful"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"gul","word2":"guld"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"syr","word2":"syr"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ful","word2":"fyl"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"syn","word2":"søn"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"fuld","word2":"fuld"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rør","word2":"rær"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"væg","word2":"vægt"}]' },
    @{ Pid = 'qCpN9NJtOIx'; L1 = 'Spanish'; Data = '[{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"fuld","word2":"fugl"},{"participantResponse":"different","correctResponse":"same","isCorrect":false,"word1":"møl","word2":"møl"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"fuld","word2":"fuld"},{"participantResponse":"different","correctResponse":"same","isCorrect":false,"word1":"røre","word2":"røre"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"skæl","word2":"skæl"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rør","word2":"rær"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ben","word2":"bøn"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"tat","word2":"tæt"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"dør","word2":"dør"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ful","word2":"fyl"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"pil","word2":"bil"},{"participantResponse":"different","correctResponse":"same","isCorrect":false,"word1":"ful","word2":"ful"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"kat","word2":"kæt"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"grå","word2":"grå"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"nø","word2":"nu"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"lidt","word2":"let"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ven","word2":"pen"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"set","word2":"sæt"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"fisk","word2":"fæsk"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"mølle","word2":"mulle"},{"participantResponse":"different","correctResponse":"same","isCorrect":false,"word1":"nø","word2":"nø"},{"participantResponse":"different","correctResponse":"same","isCorrect":false,"word1":"syr","word2":"syr"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"lys","word2":"lys"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"kat","word2":"kat"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"før","word2":"fær"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"mæt","word2":"mæt"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rød","word2":"ryd"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"løg","word2":"ly"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"rød","word2":"rød"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rød","word2":"lød"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"væg","word2":"vægt"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"syn","word2":"søn"},{"participantResponse":"different","correctResponse":"same","isCorrect":false,"word1":"ben","word2":"ben"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"syd","word2":"sød"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"sø","word2":"su"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"sinde","word2":"sende"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"lys","word2":"los"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"gul","word2":"guld"}]' },
    @{ Pid = '645a1Zg2j0U'; L1 = 'Spanish'; Data = '[{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"ven","word2":"pen"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rør","word2":"rær"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"mæt","word2":"mæt"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"set","word2":"sæt"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"mølle","word2":"mulle"},{"participantResponse":"different","correctResponse":"same","isCorrect":false,"word1":"ful","word2":"ful"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"ful","word2":"fyl"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"lys","word2":"los"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"ben","word2":"bøn"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rød","word2":"lød"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"fisk","word2":"fæsk"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"kat","word2":"kat"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"røre","word2":"røre"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"lys","word2":"lys"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"tat","word2":"tæt"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"sø","word2":"su"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"løg","word2":"ly"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"dør","word2":"dør"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"fuld","word2":"fuld"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"pil","word2":"bil"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"rød","word2":"ryd"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"syn","word2":"søn"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"fuld","word2":"fugl"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"kat","word2":"kæt"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"møl","word2":"møl"},{"participantResponse":"different","correctResponse":"same","isCorrect":false,"word1":"nø","word2":"nø"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"sinde","word2":"sende"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"væg","word2":"vægt"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"skæl","word2":"skæl"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"ben","word2":"ben"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"lidt","word2":"let"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"grå","word2":"grå"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"gul","word2":"guld"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"rød","word2":"rød"},{"participantResponse":"same","correctResponse":"different","isCorrect":false,"word1":"syd","word2":"sød"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"før","word2":"fær"},{"participantResponse":"same","correctResponse":"same","isCorrect":true,"word1":"syr","word2":"syr"},{"participantResponse":"different","correctResponse":"different","isCorrect":true,"word1":"nø","word2":"nu"}]' }
)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$appendAt = $lastRow + 1
foreach ($row in $spanishRows) {
    $ws.Rows.Item($appendAt).Insert()
    $ws.Range("A" + $appendAt).Value = $row.Pid
    $ws.Range("B" + $appendAt).Value = $row.L1
    $ws.Range("C" + $appendAt).Value = $row.Data
    $appendAt = $appendAt + 1
}
